# Dev Computer 16/05/2023 01
#
# The "Member Category" column (I) had its value changed from
# "Registered Member" to "Membership" for both data rows, and the
# selection/active cell was left on I2, the cell that was edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

$ws.Range("I2").Value = "Membership"
$ws.Range("I3").Value = "Membership"

# Leave the selection on the cell that was edited, matching where the
# author's cursor ended up after making the change.
$ws.Range("I2").Select() | Out-Null
